$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(16, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1911", 22666, 781242),
  @(17, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1910", 40000, 781242),
  @(18, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1909", 40000, 781242),
  @(19, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1908", 40000, 781242),
  @(20, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1907", 40000, 781242),
  @(21, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1906", 40000, 781242),
  @(22, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1905", 40000, 781242),
  @(23, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1904", 40000, 781242),
  @(24, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1903", 40000, 781242),
  @(25, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1902", 40000, 781242),
  @(26, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1901", 40000, 781242),
  @(27, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1812", 40000, 781242),
  @(28, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1811", 40000, 781242),
  @(29, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1810", 40000, 781242),
  @(30, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1809", 40000, 781242),
  @(31, "CC", "1047371050", "RONAL RAFAEL SALAS GUERRA", "1808", 40000, 781242),
  @(32, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1911", 22666, 1000000),
  @(33, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1910", 40000, 1000000),
  @(34, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1909", 40000, 1000000),
  @(35, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1908", 40000, 1000000),
  @(36, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1907", 40000, 1000000),
  @(37, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1906", 40000, 1000000),
  @(38, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1905", 40000, 1000000),
  @(39, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1904", 40000, 1000000),
  @(40, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1903", 40000, 1000000),
  @(41, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1902", 40000, 1000000),
  @(42, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1901", 40000, 1000000),
  @(43, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1812", 40000, 1000000),
  @(44, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1811", 40000, 1000000),
  @(45, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1810", 40000, 1000000),
  @(46, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1809", 40000, 1000000),
  @(47, "CC", "1047489473", "ELEAZAR DE JESUS ALCANTARA PEREZ", "1808", 40000, 1000000),
  @(48, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1911", 22666, 1000000),
  @(49, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1910", 40000, 1000000),
  @(50, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1909", 40000, 1000000),
  @(51, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1908", 40000, 1000000),
  @(52, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1907", 40000, 1000000),
  @(53, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1906", 40000, 1000000),
  @(54, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1905", 40000, 1000000),
  @(55, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1904", 40000, 1000000),
  @(56, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1903", 40000, 1000000),
  @(57, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1902", 40000, 1000000),
  @(58, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1901", 40000, 1000000),
  @(59, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1812", 40000, 1000000),
  @(60, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1811", 40000, 1000000),
  @(61, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1810", 40000, 1000000),
  @(62, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1809", 40000, 1000000),
  @(63, "CC", "13816726", "DAGOBERTO COLEY ESTEVEZ", "1808", 40000, 1000000)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $ws.Cells.Item($r, 6).Value2 = $row[5]
    $ws.Cells.Item($r, 7).Value2 = $row[6]
}

Write-Output "Done updating rows 16-63"
